$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13)
}

# --- 1) Relocate the "_GoBack" bookmark ------------------------------------
# Today it sits right after "...cash transaction with no deduction" at the
# end of the "Credit card processor service" bullet. Remove it there; it
# gets re-added below, inside the "Use mocks" sentence.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2) Delete the whole "Fraud System" sub-section ------------------------
# Four consecutive paragraphs: "Fraud System", "Photo system photos license
# plate entering and leaving (mock it) .", the "Application logic: ..."
# bullet, and the trailing empty paragraph that followed them.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq "Fraud System") {
        $lastPara = $d.Paragraphs.Item($i + 3)
        $deleteRange = $d.Range($p.Range.Start, $lastPara.Range.End)
        $deleteRange.Delete()
        break
    }
}

# --- 3) Split "Use mocks to mock out the following;" into two runs and ----
# drop the "_GoBack" bookmark right in between the two pieces.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-ParaText $p) -eq "Use mocks to mock out the following;") {
        $splitAt = $p.Range.Start + "Use mocks".Length
        $tailRange = $d.Range($splitAt, $p.Range.End - 1)
        $tailRange.Text = " to mock out the following;"
        $d.Bookmarks.Add("_GoBack", $d.Range($splitAt, $splitAt))
        break
    }
}
